# 141: 31/12 11:50 LP1912+6203+6173
# Appends the latest scrape batch to the three route sheets and refreshes
# the "Ultima actualizacion" / "Total filas" header cells.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "LP1912" (columns: A meta, B Hora_Scrap, C Hora_Llegada, D Linea,
#                 E Minutos, F Parada, G Fecha)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 08:50:45"
$ws1.Range("A3").Value = "Total filas: 791"

$data1 = @(
    @(777, "08:50:34", "09:02", "17X38_ROMERO",   12),
    @(778, "08:50:34", "09:03", "23_HERNANDEZ",   13),
    @(779, "08:50:34", "09:08", "16_SANTA ANA",   18),
    @(780, "08:50:34", "09:14", "11_ETCHEVERRY",  24),
    @(781, "08:50:34", "09:16", "27_EL RETIRO",   26),
    @(782, "08:50:34", "09:17", "23_HERNANDEZ",   27),
    @(783, "08:50:34", "09:21", "16_SANTA ANA",   31),
    @(784, "08:50:34", "09:27", "215_EL PELIGRO", 37),
    @(785, "08:50:34", "09:33", "23_HERNANDEZ",   43),
    @(786, "08:50:34", "09:44", "14_ABASTO",      54),
    @(787, "08:50:34", "09:51", "15_ABASTO",      61),
    @(788, "08:50:34", "09:54", "10_OLMOS",       64),
    @(789, "08:50:34", "10:03", "215C_EL PATO",   73),
    @(790, "08:50:34", "10:04", "14_ABASTO",      74),
    @(791, "08:50:34", "10:24", "11_ETCHEVERRY",  94),
    @(792, "08:50:34", "10:26", "15X38_ABASTO",   96)
)

foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = "LP1912"
    $ws1.Cells.Item($r, 7).Value = "31/12/2025"
}

# ----------------------------------------------------------------------
# Sheet "LP1912-215" (columns: A meta, B Fecha, C Hora_Scrap,
#                     D Hora_Llegada, E Linea, F Minutos, G Parada)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 08:50:45"
$ws2.Range("A3").Value = "Total filas: 58"

$data2 = @(
    @(58, "08:50:34", "09:27", "215_EL PELIGRO", 37),
    @(59, "08:50:34", "10:03", "215C_EL PATO",   73)
)

foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = "31/12/2025"
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = "LP1912"
}

# ----------------------------------------------------------------------
# Sheet "6203-6173" (columns: A meta, B Fecha, C Hora_Scrap,
#                    D Hora_Llegada, E Linea, F Minutos, G Parada)
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 08:50:45"
$ws3.Range("A3").Value = "Total filas: 93"

$ws3.Cells.Item(93, 2).Value = "31/12/2025"
$ws3.Cells.Item(93, 3).Value = "08:50:39"
$ws3.Cells.Item(93, 4).Value = "10:08"
$ws3.Cells.Item(93, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(93, 6).Value = 78
$ws3.Cells.Item(93, 7).Value = "L6203"

$ws3.Cells.Item(94, 2).Value = "31/12/2025"
$ws3.Cells.Item(94, 3).Value = "08:50:44"
$ws3.Cells.Item(94, 4).Value = "10:22"
$ws3.Cells.Item(94, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(94, 6).Value = 92
$ws3.Cells.Item(94, 7).Value = "L6173"
